# Updated cryptos list on Mon May 22 18:40:13 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# for rows 13-15 and 43-51 the ranking reshuffled so Coin (B) / Link (C) also
# move to a different row.
#
# All of these cells are stored as plain text in the workbook (prices like
# "27.020.94" or "1.005" are NOT numbers here - note the multiple "."s used
# as thousands separators, and the %-strings carry padding spaces). A plain
# `Range.Value = "..."` assignment lets Excel auto-detect/convert numeric-
# looking text into a real Number (and forcing text via NumberFormat = "@"
# mints a brand-new style record). Writing the literal as a `="..."` formula
# and immediately collapsing it to a static value with Copy + PasteSpecial
# (xlPasteValues) keeps the cell as plain text using the existing default
# style, so only the cell value changes - matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163

# Row 2
$ws.Range("D2").Formula = "=""27.020.94"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteValues)
$ws.Range("E2").Formula = "=""  -0.27%  """
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial($xlPasteValues)

# Row 3
$ws.Range("D3").Formula = "=""1.830.69"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteValues)
$ws.Range("E3").Formula = "=""  +0.39%  """
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteValues)

# Row 4
$ws.Range("E4").Formula = "=""  -0.47%  """
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial($xlPasteValues)

# Row 5
$ws.Range("D5").Formula = "=""312.30"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").Formula = "=""  +0.08%  """
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial($xlPasteValues)

# Row 6
$ws.Range("E6").Formula = "=""  -0.44%  """
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial($xlPasteValues)

# Row 7
$ws.Range("D7").Formula = "=""0.4613"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial($xlPasteValues)
$ws.Range("E7").Formula = "=""  -0.22%  """
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial($xlPasteValues)

# Row 8
$ws.Range("E8").Formula = "=""  +1.79%  """
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial($xlPasteValues)

# Row 9
$ws.Range("D9").Formula = "=""0.07347"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial($xlPasteValues)
$ws.Range("E9").Formula = "=""  +0.63%  """
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial($xlPasteValues)

# Row 10
$ws.Range("D10").Formula = "=""0.8761"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial($xlPasteValues)

# Row 11
$ws.Range("D11").Formula = "=""0.07936"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial($xlPasteValues)
$ws.Range("E11").Formula = "=""  +4.53%  """
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial($xlPasteValues)

# Row 12
$ws.Range("E12").Formula = "=""  -1.41%  """
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial($xlPasteValues)

# Row 13
$ws.Range("B13").Formula = "=""WrappedEther"""
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial($xlPasteValues)
$ws.Range("C13").Formula = "=""https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"""
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial($xlPasteValues)
$ws.Range("D13").Formula = "=""1.781.74"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial($xlPasteValues)
$ws.Range("E13").Formula = "=""  -2.75%  """
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial($xlPasteValues)

# Row 14
$ws.Range("B14").Formula = "=""Polkadot"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial($xlPasteValues)
$ws.Range("C14").Formula = "=""https://coinranking.com/coin/25W7FG7om+polkadot-dot"""
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial($xlPasteValues)
$ws.Range("D14").Formula = "=""5.341"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial($xlPasteValues)
$ws.Range("E14").Formula = "=""  -0.12%  """
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial($xlPasteValues)

# Row 15
$ws.Range("B15").Formula = "=""Chainlink"""
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial($xlPasteValues)
$ws.Range("C15").Formula = "=""https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"""
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial($xlPasteValues)
$ws.Range("D15").Formula = "=""6.562"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteValues)
$ws.Range("E15").Formula = "=""  +1.13%  """
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial($xlPasteValues)

# Row 17
$ws.Range("D17").Formula = "=""1.007"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial($xlPasteValues)
$ws.Range("E17").Formula = "=""  -0.31%  """
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial($xlPasteValues)

# Row 18
$ws.Range("D18").Formula = "=""0.000008903"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)
$ws.Range("E18").Formula = "=""  +3.17%  """
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial($xlPasteValues)

# Row 19
$ws.Range("E19").Formula = "=""  -0.49%  """
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial($xlPasteValues)

# Row 20
$ws.Range("D20").Formula = "=""14.81"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial($xlPasteValues)
$ws.Range("E20").Formula = "=""  +2.47%  """
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial($xlPasteValues)

# Row 21
$ws.Range("D21").Formula = "=""27.041.16"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteValues)
$ws.Range("E21").Formula = "=""  -1.17%  """
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial($xlPasteValues)

# Row 22
$ws.Range("D22").Formula = "=""5.119"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial($xlPasteValues)
$ws.Range("E22").Formula = "=""  -1.69%  """
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial($xlPasteValues)

# Row 23
$ws.Range("D23").Formula = "=""10.57"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial($xlPasteValues)
$ws.Range("E23").Formula = "=""  +0.03%  """
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial($xlPasteValues)

# Row 24
$ws.Range("D24").Formula = "=""2.047.58"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial($xlPasteValues)
$ws.Range("E24").Formula = "=""  -2.03%  """
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial($xlPasteValues)

# Row 25
$ws.Range("D25").Formula = "=""153.34"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteValues)
$ws.Range("E25").Formula = "=""  +0.85%  """
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial($xlPasteValues)

# Row 26
$ws.Range("E26").Formula = "=""  -1.57%  """
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial($xlPasteValues)

# Row 27
$ws.Range("D27").Formula = "=""18.44"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteValues)
$ws.Range("E27").Formula = "=""  +1.05%  """
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial($xlPasteValues)

# Row 28
$ws.Range("E28").Formula = "=""  -1.96%  """
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial($xlPasteValues)

# Row 29
$ws.Range("D29").Formula = "=""5.158"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial($xlPasteValues)
$ws.Range("E29").Formula = "=""  +1.49%  """
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial($xlPasteValues)

# Row 30
$ws.Range("D30").Formula = "=""115.65"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial($xlPasteValues)

# Row 31
$ws.Range("D31").Formula = "=""0.08911"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial($xlPasteValues)
$ws.Range("E31").Formula = "=""  +0.07%  """
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial($xlPasteValues)

# Row 32
$ws.Range("D32").Formula = "=""2.963"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial($xlPasteValues)
$ws.Range("E32").Formula = "=""  +0.15%  """
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial($xlPasteValues)

# Row 33
$ws.Range("D33").Formula = "=""0.7327"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial($xlPasteValues)

# Row 34
$ws.Range("D34").Formula = "=""4.438"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial($xlPasteValues)
$ws.Range("E34").Formula = "=""  -0.39%  """
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial($xlPasteValues)

# Row 35
$ws.Range("E35").Formula = "=""  -0.47%  """
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial($xlPasteValues)

# Row 36
$ws.Range("D36").Formula = "=""2.485"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$ws.Range("E36").Formula = "=""  +0.35%  """
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial($xlPasteValues)

# Row 37
$ws.Range("E37").Formula = "=""  +1.99%  """
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial($xlPasteValues)

# Row 38
$ws.Range("D38").Formula = "=""1.072"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial($xlPasteValues)
$ws.Range("E38").Formula = "=""  -0.14%  """
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial($xlPasteValues)

# Row 39
$ws.Range("D39").Formula = "=""0.05240"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial($xlPasteValues)
$ws.Range("E39").Formula = "=""  -0.20%  """
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial($xlPasteValues)

# Row 40
$ws.Range("D40").Formula = "=""2.942"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial($xlPasteValues)
$ws.Range("E40").Formula = "=""  +0.57%  """
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial($xlPasteValues)

# Row 41
$ws.Range("D41").Formula = "=""7.120"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial($xlPasteValues)
$ws.Range("E41").Formula = "=""  -0.21%  """
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial($xlPasteValues)

# Row 42
$ws.Range("D42").Formula = "=""0.5161"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial($xlPasteValues)
$ws.Range("E42").Formula = "=""  -0.72%  """
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial($xlPasteValues)

# Row 43
$ws.Range("B43").Formula = "=""Algorand"""
$ws.Range("B43").Copy()
$ws.Range("B43").PasteSpecial($xlPasteValues)
$ws.Range("C43").Formula = "=""https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"""
$ws.Range("C43").Copy()
$ws.Range("C43").PasteSpecial($xlPasteValues)
$ws.Range("D43").Formula = "=""0.1627"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$ws.Range("E43").Formula = "=""  -0.23%  """
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial($xlPasteValues)

# Row 44
$ws.Range("B44").Formula = "=""Aptos"""
$ws.Range("B44").Copy()
$ws.Range("B44").PasteSpecial($xlPasteValues)
$ws.Range("C44").Formula = "=""https://coinranking.com/coin/HGYj5JCv5+aptos-apt"""
$ws.Range("C44").Copy()
$ws.Range("C44").PasteSpecial($xlPasteValues)
$ws.Range("D44").Formula = "=""8.212"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial($xlPasteValues)
$ws.Range("E44").Formula = "=""  -0.64%  """
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial($xlPasteValues)

# Row 45
$ws.Range("B45").Formula = "=""Decentraland"""
$ws.Range("B45").Copy()
$ws.Range("B45").PasteSpecial($xlPasteValues)
$ws.Range("C45").Formula = "=""https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"""
$ws.Range("C45").Copy()
$ws.Range("C45").PasteSpecial($xlPasteValues)
$ws.Range("D45").Formula = "=""0.4846"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial($xlPasteValues)
$ws.Range("E45").Formula = "=""  -0.39%  """
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial($xlPasteValues)

# Row 46
$ws.Range("B46").Formula = "=""EnergySwap"""
$ws.Range("B46").Copy()
$ws.Range("B46").PasteSpecial($xlPasteValues)
$ws.Range("C46").Formula = "=""https://coinranking.com/coin/SbWqqTui-+energyswap-ens"""
$ws.Range("C46").Copy()
$ws.Range("C46").PasteSpecial($xlPasteValues)
$ws.Range("D46").Formula = "=""10.23"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial($xlPasteValues)
$ws.Range("E46").Formula = "=""  +0.88%  """
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial($xlPasteValues)

# Row 47
$ws.Range("B47").Formula = "=""PaxDollar"""
$ws.Range("B47").Copy()
$ws.Range("B47").PasteSpecial($xlPasteValues)
$ws.Range("C47").Formula = "=""https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"""
$ws.Range("C47").Copy()
$ws.Range("C47").PasteSpecial($xlPasteValues)
$ws.Range("D47").Formula = "=""1.005"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial($xlPasteValues)
$ws.Range("E47").Formula = "=""  -0.45%  """
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial($xlPasteValues)

# Row 48
$ws.Range("B48").Formula = "=""Quant"""
$ws.Range("B48").Copy()
$ws.Range("B48").PasteSpecial($xlPasteValues)
$ws.Range("C48").Formula = "=""https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"""
$ws.Range("C48").Copy()
$ws.Range("C48").PasteSpecial($xlPasteValues)
$ws.Range("D48").Formula = "=""102.41"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial($xlPasteValues)
$ws.Range("E48").Formula = "=""  -1.18%  """
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial($xlPasteValues)

# Row 49
$ws.Range("B49").Formula = "=""NEARProtocol"""
$ws.Range("B49").Copy()
$ws.Range("B49").PasteSpecial($xlPasteValues)
$ws.Range("C49").Formula = "=""https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"""
$ws.Range("C49").Copy()
$ws.Range("C49").PasteSpecial($xlPasteValues)
$ws.Range("D49").Formula = "=""1.634"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial($xlPasteValues)
$ws.Range("E49").Formula = "=""  +0.20%  """
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial($xlPasteValues)

# Row 50
$ws.Range("B50").Formula = "=""Cronos"""
$ws.Range("B50").Copy()
$ws.Range("B50").PasteSpecial($xlPasteValues)
$ws.Range("C50").Formula = "=""https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"""
$ws.Range("C50").Copy()
$ws.Range("C50").PasteSpecial($xlPasteValues)
$ws.Range("D50").Formula = "=""0.06203"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$ws.Range("E50").Formula = "=""  -0.87%  """
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial($xlPasteValues)

# Row 51
$ws.Range("B51").Formula = "=""Aave"""
$ws.Range("B51").Copy()
$ws.Range("B51").PasteSpecial($xlPasteValues)
$ws.Range("C51").Formula = "=""https://coinranking.com/coin/ixgUfzmLR+aave-aave"""
$ws.Range("C51").Copy()
$ws.Range("C51").PasteSpecial($xlPasteValues)
$ws.Range("D51").Formula = "=""64.65"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial($xlPasteValues)
$ws.Range("E51").Formula = "=""  +0.20%  """
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial($xlPasteValues)

$excel.CutCopyMode = $false

